# "Add BDD Others Files" — drop the "-name" option from the `find`
# command-line examples (find /usr/bin -name *td  ->  find /usr/bin *td),
# and turn the "||" between the two find commands into "&&".

$d = $word.ActiveDocument

# Deletes the first occurrence of $substr inside paragraph $paraIndex.
function Remove-Substring($doc, $paraIndex, $substr) {
    $p = $doc.Paragraphs($paraIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($substr)
    if ($idx -lt 0) { return }
    $start = $p.Range.Start + $idx
    $end = $start + $substr.Length
    $sub = $doc.Range($start, $end)
    $sub.Text = ""
}

# Replaces the first occurrence of $substr inside paragraph $paraIndex with
# $newtext, forcing the new text to land in its own run (distinct from the
# runs that surround it) by nudging a character-formatting property on it
# and then reverting the nudge — mirrors what happens when a run gets
# split in two by typing in the middle of it.
function Replace-WithSplit($doc, $paraIndex, $substr, $newtext) {
    $p = $doc.Paragraphs($paraIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($substr)
    if ($idx -lt 0) { return }
    $start = $p.Range.Start + $idx
    $end = $start + $substr.Length
    $sub = $doc.Range($start, $end)
    $sub.Text = ""
    $collapsed = $doc.Range($start, $start)
    $collapsed.InsertAfter($newtext)
    $newRange = $doc.Range($start, $start + $newtext.Length)
    $newRange.Font.Bold = $true
    $newRange.Font.Bold = $false
}

# Locate (by content, not a hard-coded index) every paragraph that still
# contains the "-name" find(1) option and fix it up.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs($i).Range.Text

    if ($text.IndexOf("-name") -lt 0) {
        continue
    }

    if ($text.IndexOf("||") -ge 0) {
        # "find /usr/bin -name *td || find /usr/bin -name *cd"
        #   -> "find /usr/bin *td && find /usr/bin *cd"
        Remove-Substring $d $i "-name "
        Replace-WithSplit $d $i "||" "&&"
        Remove-Substring $d $i "-name "
    } else {
        # "find /usr/bin -name *td"  -> "find /usr/bin *td"
        # "find /usr/bin -name *td*" -> "find /usr/bin *td*"
        Remove-Substring $d $i "-name "
    }
}
